$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 7 (old row 7 and
# everything below shifts down by two rows: old row 7 -> new row 9, ...,
# old row 117 -> new row 119). This matches the dimension growing from
# A1:T117 to A1:T119.
$ws.Range("A7:T8").EntireRow.Insert()

# Fill in the first new row (new row 7) with its data.
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Macroferia Regional de Talca"
$ws.Range("C7").Value = "Maule"
$ws.Range("D7").Value = 44963
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100101
$ws.Range("H7").Value = "Berries"
$ws.Range("I7").Value = 100101001
$ws.Range("J7").Value = "Arándano (blue)"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 150
$ws.Range("N7").Value = 3000
$ws.Range("O7").Value = 3000
$ws.Range("P7").Value = 3000
$ws.Range("Q7").Value = "`$/bandeja 2 kilos"
$ws.Range("R7").Value = "Provincia de Curicó"
$ws.Range("S7").Value = 1500
$ws.Range("T7").Value = 2

# Fill in the second new row (new row 8) with its data.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = 44963
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100101
$ws.Range("H8").Value = "Berries"
$ws.Range("I8").Value = 100101001
$ws.Range("J8").Value = "Arándano (blue)"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 2500
$ws.Range("O8").Value = 2500
$ws.Range("P8").Value = 2500
$ws.Range("Q8").Value = "`$/bandeja 2 kilos"
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value = 1250
$ws.Range("T8").Value = 2
